$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($sheet, $addr, $text) {
    $range = $sheet.Range($addr)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

Set-CellText $ws "D2" "326.69"
Set-CellText $ws "E2" "-1.31%"
Set-CellText $ws "D3" "44.25"
Set-CellText $ws "E3" "0.17%"
Set-CellText $ws "D4" "5.236"
Set-CellText $ws "E4" "-5.50%"
Set-CellText $ws "D5" "0.08336"
Set-CellText $ws "E5" "2.29%"
Set-CellText $ws "B6" "FTXToken"
Set-CellText $ws "C6" "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-CellText $ws "D6" "1.942"
Set-CellText $ws "E6" "-5.59%"
Set-CellText $ws "B7" "MXToken"
Set-CellText $ws "C7" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-CellText $ws "D7" "0.9700"
Set-CellText $ws "E7" "-0.70%"
Set-CellText $ws "B8" "BTSEToken"
Set-CellText $ws "C8" "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-CellText $ws "D8" "2.541"
Set-CellText $ws "E8" "-2.72%"
Set-CellText $ws "B9" "LiechtensteinCryptoassetsExchange"
Set-CellText $ws "C9" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-CellText $ws "D9" "0.1127"
Set-CellText $ws "E9" "1.76%"
Set-CellText $ws "B10" "WazirX"
Set-CellText $ws "C10" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-CellText $ws "D10" "0.1888"
Set-CellText $ws "E10" "-0.28%"
Set-CellText $ws "B11" "MandalaExchangeToken"
Set-CellText $ws "C11" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-CellText $ws "D11" "0.09715"
Set-CellText $ws "E11" "-2.95%"
Set-CellText $ws "B12" "BitrueCoin"
Set-CellText $ws "C12" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-CellText $ws "D12" "0.04608"
Set-CellText $ws "E12" "-2.34%"
Set-CellText $ws "B13" "BitMartToken"
Set-CellText $ws "C13" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-CellText $ws "D13" "0.1059"
Set-CellText $ws "E13" "0.11%"
Set-CellText $ws "B14" "BitForexToken"
Set-CellText $ws "C14" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-CellText $ws "D14" "0.001280"
Set-CellText $ws "E14" "1.43%"
Set-CellText $ws "B15" "TigerCash"
Set-CellText $ws "C15" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-CellText $ws "D15" "0.005892"
Set-CellText $ws "E15" "-1.18%"
Set-CellText $ws "B16" "LEO"
Set-CellText $ws "C16" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-CellText $ws "D16" "3.403"
Set-CellText $ws "E16" "1.81%"
Set-CellText $ws "B17" "GateToken"
Set-CellText $ws "C17" "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-CellText $ws "D17" "4.429"
Set-CellText $ws "E17" "0.00%"
Set-CellText $ws "E18" "-0.11%"
Set-CellText $ws "D19" "8.601"
Set-CellText $ws "E19" "-15.66%"
Set-CellText $ws "D20" "0.1372"
Set-CellText $ws "E20" "-1.27%"
Set-CellText $ws "D21" "0.2579"
Set-CellText $ws "E21" "0.41%"
Set-CellText $ws "D22" "0.04155"
Set-CellText $ws "E22" "1.39%"
Set-CellText $ws "E23" "-5.45%"
Set-CellText $ws "D24" "0.004427"
Set-CellText $ws "E24" "0.99%"
Set-CellText $ws "D25" "0.0001303"
Set-CellText $ws "E25" "1.92%"
Set-CellText $ws "D26" "0.0002985"
Set-CellText $ws "E26" "-20.05%"
Set-CellText $ws "D38" "0.02692"
Set-CellText $ws "E38" "0.11%"
Set-CellText $ws "D39" "0.05563"
Set-CellText $ws "E39" "-1.54%"
Set-CellText $ws "D40" "0.007833"
Set-CellText $ws "E40" "2.66%"
Set-CellText $ws "D41" "0.1409"
Set-CellText $ws "E41" "-0.82%"
Set-CellText $ws "D42" "0.007310"
Set-CellText $ws "E42" "-3.31%"
Set-CellText $ws "D43" "0.002114"
Set-CellText $ws "E43" "8.08%"
Set-CellText $ws "D44" "0.007855"
Set-CellText $ws "E44" "-5.32%"
Set-CellText $ws "D45" "0.3498"
Set-CellText $ws "D46" "0.00006874"
Set-CellText $ws "E46" "-2.01%"
Set-CellText $ws "E47" "0.36%"
Set-CellText $ws "D48" "0.003491"
Set-CellText $ws "E48" "-0.95%"
Set-CellText $ws "D49" "0.003538"
Set-CellText $ws "E49" "40.60%"
Set-CellText $ws "E50" "0.36%"
Set-CellText $ws "D51" "0.0002004"
Set-CellText $ws "E51" "0.36%"
